# Apply the "Add files via upload" edit:
#  - bump every cached "datetime2" auto-date field (slide master + every
#    slide layout) from "Sunday, April 3, 2022" to "Monday, March 6, 2023"
#  - update the subtitle on slide 1, dropping the trailing student id
#    ("DIAMANTIS RAFAIL PAPADAM - 2017030044" -> "DIAMANTIS RAFAIL PAPADAM")

$p = $ppt.ActivePresentation

$oldDateText = "Sunday, April 3, 2022"
$newDateText = "Monday, March 6, 2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDateText) {
                    $shp.TextFrame.TextRange.Text = $newDateText
                }
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every custom (slide) layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# Slide 1: subtitle run "DIAMANTIS RAFAIL PAPADAM - 2017030044" -> "DIAMANTIS RAFAIL PAPADAM"
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(3)
$oldName = "DIAMANTIS RAFAIL PAPADAM - 2017030044"
$newName = "DIAMANTIS RAFAIL PAPADAM"
if ($subtitle.TextFrame.TextRange.Text -like "$oldName*") {
    $subtitle.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = $newName
}
